$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the AF (Phases) column values per the diff.

$ws.Range("AF2").Value = "DICHOTIC_PRE"
$ws.Range("AF3").Value = "DIGIT_PRE"
$ws.Range("AF4").Value = "DICHOTIC_POST"
$ws.Range("AF5").Value = "Dichotic_before_after"
$ws.Range("AF6").Value = "DICHOTIC_PRE"
$ws.Range("AF7").Value = "DIGIT_PRE"
$ws.Range("AF8").Value = "Dichotic_before_after"
$ws.Range("AF9").Value = "DIGIT_PRE"
$ws.Range("AF10").Value = "DICHOTIC_PRE"
$ws.Range("AF11").Value = "DIGIT_PRE"
$ws.Range("AF12").Value = "DICHOTIC_POST"
$ws.Range("AF13").Value = "DICHOTIC_PRE"
$ws.Range("AF14").Value = "Digit_before_after"
$ws.Range("AF15").Value = "Dichotic_before_after"
$ws.Range("AF16").Value = "DICHOTIC_POST"
$ws.Range("AF17").Value = "DIGIT_POST"
$ws.Range("AF18").Value = "Digit_before_after"
$ws.Range("AF19").Value = "DIGIT_POST"
$ws.Range("AF20").Value = "Digit_before_after"
$ws.Range("AF21").Value = "DICHOTIC_PRE"
$ws.Range("AF22").Value = "DICHOTIC_POST"
$ws.Range("AF23").Value = "Digit_before_after"
$ws.Range("AF24").Value = "DIGIT_POST"
$ws.Range("AF25").Value = "Dichotic_before_after"
$ws.Range("AF26").Value = "DIGIT_POST"
$ws.Range("AF27").Value = "DIGIT_PRE"
$ws.Range("AF28").Value = "Dichotic_before_after"
$ws.Range("AF29").Value = "DIGIT_POST"
$ws.Range("AF30").Value = "DICHOTIC_POST"
$ws.Range("AF31").Value = "DICHOTIC_PRE"
$ws.Range("AF32").Value = "DIGIT_PRE"
$ws.Range("AF33").Value = "Digit_before_after"
$ws.Range("AF34").Value = "DIGIT_POST"
$ws.Range("AF35").Value = "DIGIT_POST"
$ws.Range("AF36").Value = "Dichotic_before_after"
$ws.Range("AF37").Value = "DICHOTIC_POST"
$ws.Range("AF38").Value = "Digit_before_after"
$ws.Range("AF39").Value = "DIGIT_PRE"
$ws.Range("AF40").Value = "DICHOTIC_PRE"
$ws.Range("AF41").Value = "DICHOTIC_POST"
$ws.Range("AF42").Value = "DIGIT_PRE"
$ws.Range("AF43").Value = "DIGIT_PRE"
$ws.Range("AF44").Value = "DIGIT_POST"
$ws.Range("AF45").Value = "DIGIT_PRE"
$ws.Range("AF46").Value = "DIGIT_PRE"
$ws.Range("AF47").Value = "Digit_before_after"
$ws.Range("AF48").Value = "Dichotic_before_after"
$ws.Range("AF49").Value = "DICHOTIC_PRE"
$ws.Range("AF50").Value = "DICHOTIC_PRE"
$ws.Range("AF51").Value = "Digit_before_after"
$ws.Range("AF52").Value = "DIGIT_POST"
$ws.Range("AF53").Value = "Digit_before_after"
$ws.Range("AF54").Value = "DIGIT_POST"
$ws.Range("AF55").Value = "DIGIT_PRE"
$ws.Range("AF56").Value = "DIGIT_POST"
$ws.Range("AF57").Value = "DICHOTIC_POST"
$ws.Range("AF58").Value = "Dichotic_before_after"
$ws.Range("AF59").Value = "DICHOTIC_PRE"
$ws.Range("AF60").Value = "Dichotic_before_after"
$ws.Range("AF61").Value = "DICHOTIC_PRE"
$ws.Range("AF62").Value = "DICHOTIC_POST"
$ws.Range("AF63").Value = "DICHOTIC_POST"
$ws.Range("AF64").Value = "DICHOTIC_POST"
$ws.Range("AF65").Value = "Dichotic_before_after"
$ws.Range("AF66").Value = "DIGIT_PRE"
$ws.Range("AF67").Value = "DICHOTIC_PRE"
$ws.Range("AF68").Value = "DIGIT_POST"
$ws.Range("AF69").Value = "DIGIT_POST"
$ws.Range("AF70").Value = "Digit_before_after"
$ws.Range("AF71").Value = "DICHOTIC_POST"
$ws.Range("AF72").Value = "DIGIT_PRE"
$ws.Range("AF73").Value = "DIGIT_POST"
$ws.Range("AF74").Value = "DICHOTIC_POST"
$ws.Range("AF75").Value = "Digit_before_after"
$ws.Range("AF76").Value = "DICHOTIC_PRE"
$ws.Range("AF77").Value = "Dichotic_before_after"
$ws.Range("AF78").Value = "DICHOTIC_PRE"
$ws.Range("AF79").Value = "Dichotic_before_after"
$ws.Range("AF80").Value = "DICHOTIC_POST"
$ws.Range("AF81").Value = "Digit_before_after"
